$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Newly started topics today (stack, queue, prefix sum, set & map) - log date
$ws1.Range("E2").Value = 45665
$ws1.Range("E3").Value = 45665
$ws1.Range("E4").Value = 45665
$ws1.Range("E5").Value = 45665

# dp: logged 4 exercise attempts
$ws1.Range("C10").Value = 4

# bfs / dfs: date moved forward a day
$ws1.Range("E15").Value = 45665
$ws1.Range("E16").Value = 45665

# trie: clear its date entirely (fully removes the cell, not just its value)
$ws1.Range("E17").Clear()

# geometry row: fill in attempt/exercise counters
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 0

# new topic row: intermediate dp
$ws1.Range("A21").Value = "intermediate dp"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 0

# Match the center-aligned look used by the rest of column B/C
$ws1.Range("B20:C21").HorizontalAlignment = -4108

# Selections as left by the author
[void]$ws1.Range("E3").Select()
[void]$ws2.Range("B4").Select()
[void]$ws1.Activate()
